$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr1 = New-Object "object[,]" 1,7
$arr2 = New-Object "object[,]" 1,4

# Row 2
$arr1[0,0] = 0.1532258774662978
$arr1[0,1] = 0.05151343658369001
$arr1[0,2] = 0.1081177584050153
$arr1[0,3] = 1.262987421358872
$arr1[0,4] = 1.141031977449174
$arr1[0,5] = 1.105117920401526
$arr1[0,6] = 1.22614164198621
$ws.Range("C2:I2").Value = $arr1
$arr2[0,0] = 1.299466841442097
$arr2[0,1] = 0.141449036860628
$arr2[0,2] = 0.4094407034417884
$arr2[0,3] = 1.485015210462645
$ws.Range("K2:N2").Value = $arr2

# Row 3
$arr1[0,0] = 0.1517420198934047
$arr1[0,1] = 0.05127568612486044
$arr1[0,2] = 0.1079552272702813
$arr1[0,3] = 1.256789296655995
$arr1[0,4] = 1.136231794590458
$arr1[0,5] = 1.108650299514665
$arr1[0,6] = 1.22422393723113
$ws.Range("C3:I3").Value = $arr1
$arr2[0,0] = 1.1825441993781
$arr2[0,1] = 0.14187178329313
$arr2[0,2] = 0.386067738871418
$arr2[0,3] = 1.501590782394082
$ws.Range("K3:N3").Value = $arr2

# Row 4
$arr1[0,0] = 0.1508882270461385
$arr1[0,1] = 0.05113230160958793
$arr1[0,2] = 0.1079021266124656
$arr1[0,3] = 1.253823749281494
$arr1[0,4] = 1.13410451246537
$arr1[0,5] = 1.111406002069444
$arr1[0,6] = 1.22380469520931
$ws.Range("C4:I4").Value = $arr1
$arr2[0,0] = 1.111061978348999
$arr2[0,1] = 0.1421971039189707
$arr2[0,2] = 0.3718782068275743
$arr2[0,3] = 1.512342379987256
$ws.Range("K4:N4").Value = $arr2

# Row 5
$arr1[0,0] = 0.1505547509712244
$arr1[0,1] = 0.05107453497634395
$arr1[0,2] = 0.1078922392036041
$arr1[0,3] = 1.252826149171391
$arr1[0,4] = 1.133443329246717
$arr1[0,5] = 1.112676423478987
$arr1[0,6] = 1.223824249077026
$ws.Range("C5:I5").Value = $arr1
$arr2[0,0] = 1.082010927054284
$arr2[0,1] = 0.1423462069604966
$arr2[0,2] = 0.3661366662434844
$arr2[0,3] = 1.516868108061498
$ws.Range("K5:N5").Value = $arr2

# Row 6
$arr1[0,0] = 0.1505002516118239
$arr1[0,1] = 0.05106498328872178
$arr1[0,2] = 0.1078913074673533
$arr1[0,3] = 1.252673224514645
$arr1[0,4] = 1.133345949371531
$arr1[0,5] = 1.112896278523593
$arr1[0,6] = 1.223838987840978
$ws.Range("C6:I6").Value = $arr1
$arr2[0,0] = 1.077191792037297
$arr2[0,1] = 0.1423719638670846
$arr2[0,2] = 0.3651857580053672
$arr2[0,3] = 1.517628321005329
$ws.Range("K6:N6").Value = $arr2

# Row 7
$arr1[0,0] = 0.1508836711003809
$arr1[0,1] = 0.05113151984740227
$arr1[0,2] = 0.1079019456723742
$arr1[0,3] = 1.253809441930798
$arr1[0,4] = 1.134094763314025
$arr1[0,5] = 1.111422538543152
$arr1[0,6] = 1.223804188360013
$ws.Range("C7:I7").Value = $arr1
$arr2[0,0] = 1.110669866535943
$arr2[0,1] = 0.1421990478383961
$arr2[0,2] = 0.3718006089258239
$arr2[0,3] = 1.512402831019759
$ws.Range("K7:N7").Value = $arr2

# Row 8
$arr1[0,0] = 0.15270237502844
$arr1[0,1] = 0.05143092811203687
$arr1[0,2] = 0.1080520319164577
$arr1[0,3] = 1.260675709614247
$arr1[0,4] = 1.139206362760234
$arr1[0,5] = 1.106214025874422
$arr1[0,6] = 1.225322880414822
$ws.Range("C8:I8").Value = $arr1
$arr2[0,0] = 1.25908832976873
$arr2[0,1] = 0.1415811515282179
$arr2[0,2] = 0.4013482739233467
$arr2[0,3] = 1.490611304349073
$ws.Range("K8:N8").Value = $arr2

# Row 9
$arr1[0,0] = 0.1567219763817178
$arr1[0,1] = 0.05203821347408066
$arr1[0,2] = 0.1087165527450757
$arr1[0,3] = 1.280825922865972
$arr1[0,4] = 1.155762115898469
$arr1[0,5] = 1.100662070864004
$arr1[0,6] = 1.234331333427598
$ws.Range("C9:I9").Value = $arr1
$arr2[0,0] = 1.552562581090456
$arr2[0,1] = 0.1408913763324691
$arr2[0,2] = 0.4605684357832018
$arr2[0,3] = 1.452433536957447
$ws.Range("K9:N9").Value = $arr2

# Row 10
$arr1[0,0] = 0.1599498425941164
$arr1[0,1] = 0.05249613355148952
$arr1[0,2] = 0.1094302749999443
$arr1[0,3] = 1.299736733442572
$arr1[0,4] = 1.171945925451979
$arr1[0,5] = 1.099434697772267
$arr1[0,6] = 1.244648432112768
$ws.Range("C10:I10").Value = $arr1
$arr2[0,0] = 1.769648325423191
$arr2[0,1] = 0.1407032242779351
$arr2[0,2] = 0.5048549064961563
$arr2[0,3] = 1.427160781253363
$ws.Range("K10:N10").Value = $arr2

# Row 11
$arr1[0,0] = 0.1614776261225188
$arr1[0,1] = 0.05270689806604523
$arr1[0,2] = 0.1098039084767777
$arr1[0,3] = 1.309238339636011
$arr1[0,4] = 1.180189815872268
$arr1[0,5] = 1.099497752449167
$arr1[0,6] = 1.250150101924518
$ws.Range("C11:I11").Value = $arr1
$arr2[0,0] = 1.868725205440683
$arr2[0,1] = 0.1406869289434596
$arr2[0,2] = 0.5251708601221594
$arr2[0,3] = 1.416266400726386
$ws.Range("K11:N11").Value = $arr2

# Row 12
$arr1[0,0] = 0.1620646669828432
$arr1[0,1] = 0.05278705240981907
$arr1[0,2] = 0.1099524283879809
$arr1[0,3] = 1.312966113734348
$arr1[0,4] = 1.183438991212654
$arr1[0,5] = 1.09961114612031
$arr1[0,6] = 1.252350043656264
$ws.Range("C12:I12").Value = $arr1
$arr2[0,0] = 1.906289026409922
$arr2[0,1] = 0.1406907300655504
$arr2[0,2] = 0.532888301224844
$arr2[0,3] = 1.412227634621338
$ws.Range("K12:N12").Value = $arr2

# Row 13
$arr1[0,0] = 0.1619378597304149
$arr1[0,1] = 0.05276977467504551
$arr1[0,2] = 0.109920129282159
$arr1[0,3] = 1.312157494279191
$arr1[0,4] = 1.182733547876779
$arr1[0,5] = 1.099582740922841
$arr1[0,6] = 1.251871056589238
$ws.Range("C13:I13").Value = $arr1
$arr2[0,0] = 1.898196974496386
$arr2[0,1] = 0.1406894678109865
$arr2[0,2] = 0.5312251375333972
$arr2[0,3] = 1.413093597747711
$ws.Range("K13:N13").Value = $arr2

# Row 14
$arr1[0,0] = 0.1615257521586102
$arr1[0,1] = 0.05271348560973976
$arr1[0,2] = 0.1098159863866464
$arr1[0,3] = 1.30954242344967
$arr1[0,4] = 1.180454570938593
$arr1[0,5] = 1.099505286448959
$arr1[0,6] = 1.250328754108352
$ws.Range("C14:I14").Value = $arr1
$arr2[0,0] = 1.871814696039394
$arr2[0,1] = 0.1406870417725514
$arr2[0,2] = 0.5258052941081246
$arr2[0,3] = 1.415932390908509
$ws.Range("K14:N14").Value = $arr2

# Row 15
$arr1[0,0] = 0.1612744304292022
$arr1[0,1] = 0.05267905121232275
$arr1[0,2] = 0.1097531115576302
$arr1[0,3] = 1.307957524073046
$arr1[0,4] = 1.179075239696601
$arr1[0,5] = 1.099469505807718
$arr1[0,6] = 1.249399241322308
$ws.Range("C15:I15").Value = $arr1
$arr2[0,0] = 1.855660693214929
$arr2[0,1] = 0.140686854593838
$arr2[0,2] = 0.522488633454941
$arr2[0,3] = 1.417682527167337
$ws.Range("K15:N15").Value = $arr2

# Row 16
$arr1[0,0] = 0.1598511903491726
$arr1[0,1] = 0.05248240807723903
$arr1[0,2] = 0.1094068417217713
$arr1[0,3] = 1.299133909332198
$arr1[0,4] = 1.17142496232924
$arr1[0,5] = 1.099443094701115
$arr1[0,6] = 1.244305180201053
$ws.Range("C16:I16").Value = $arr1
$arr2[0,0] = 1.763179845377294
$arr2[0,1] = 0.1407056842331968
$arr2[0,2] = 0.5035306120780376
$arr2[0,3] = 1.427884882127664
$ws.Range("K16:N16").Value = $arr2

# Row 17
$arr1[0,0] = 0.1589932660782409
$arr1[0,1] = 0.05236239493830652
$arr1[0,2] = 0.1092069505645412
$arr1[0,3] = 1.293951477681702
$arr1[0,4] = 1.16695803405959
$arr1[0,5] = 1.099586159629411
$arr1[0,6] = 1.241387420043196
$ws.Range("C17:I17").Value = $arr1
$arr2[0,0] = 1.70652800701771
$arr2[0,1] = 0.140734988973243
$arr2[0,2] = 0.491943838831844
$arr2[0,3] = 1.434298024978268
$ws.Range("K17:N17").Value = $arr2

# Row 18
$arr1[0,0] = 0.1585054051944326
$arr1[0,1] = 0.05229359818689261
$arr1[0,2] = 0.1090965869926563
$arr1[0,3] = 1.291055269633205
$arr1[0,4] = 1.16447173815466
$arr1[0,5] = 1.099726925468275
$arr1[0,6] = 1.239785268905905
$ws.Range("C18:I18").Value = $arr1
$arr2[0,0] = 1.673973872811644
$arr2[0,1] = 0.1407583663097043
$arr2[0,2] = 0.4852954283627966
$arr2[0,3] = 1.43804338302251
$ws.Range("K18:N18").Value = $arr2

# Row 19
$arr1[0,0] = 0.1583411859852788
$arr1[0,1] = 0.0522703448967512
$arr1[0,2] = 0.1090600115206861
$arr1[0,3] = 1.290089178047552
$arr1[0,4] = 1.163644152081019
$arr1[0,5] = 1.099784625276897
$arr1[0,6] = 1.239255862500642
$ws.Range("C19:I19").Value = $arr1
$arr2[0,0] = 1.662956870975052
$arr2[0,1] = 0.1407674014241778
$arr2[0,2] = 0.4830471437703707
$arr2[0,3] = 1.439321230593123
$ws.Range("K19:N19").Value = $arr2

# Row 20
$arr1[0,0] = 0.1590840148728034
$arr1[0,1] = 0.05237514664506548
$arr1[0,2] = 0.1092277524049479
$arr1[0,3] = 1.294494398423694
$arr1[0,4] = 1.167424955546991
$arr1[0,5] = 1.099564876637444
$arr1[0,6] = 1.24169014547752
$ws.Range("C20:I20").Value = $arr1
$arr2[0,0] = 1.712555541775316
$arr2[0,1] = 0.1407311943776435
$arr2[0,2] = 0.493175616573204
$arr2[0,3] = 1.43360946737954
$ws.Range("K20:N20").Value = $arr2

# Row 21
$arr1[0,0] = 0.1616465677135466
$arr1[0,1] = 0.05273000986364451
$arr1[0,2] = 0.1098463848604716
$arr1[0,3] = 1.310307008930351
$arr1[0,4] = 1.181120499569289
$arr1[0,5] = 1.099525605968495
$arr1[0,6] = 1.250778599126718
$ws.Range("C21:I21").Value = $arr1
$arr2[0,0] = 1.879562581944697
$arr2[0,1] = 0.1406874836684651
$arr2[0,2] = 0.5273965773906895
$arr2[0,3] = 1.415096214829717
$ws.Range("K21:N21").Value = $arr2

# Row 22
$arr1[0,0] = 0.1633708838963202
$arr1[0,1] = 0.05296392450733478
$arr1[0,2] = 0.1102916876333282
$arr1[0,3] = 1.32139779361458
$arr1[0,4] = 1.190814149516086
$arr1[0,5] = 1.10002178554717
$arr1[0,6] = 1.257398097081349
$ws.Range("C22:I22").Value = $arr1
$arr2[0,0] = 1.98897643798847
$arr2[0,1] = 0.1407170430192153
$arr2[0,2] = 0.5499031005777795
$arr2[0,3] = 1.403502128780033
$ws.Range("K22:N22").Value = $arr2

# Row 23
$arr1[0,0] = 0.1624460644083285
$arr1[0,1] = 0.05283890110720435
$arr1[0,2] = 0.110050272786669
$arr1[0,3] = 1.315409079803146
$arr1[0,4] = 1.185572301732748
$arr1[0,5] = 1.099709162645553
$arr1[0,6] = 1.253802845052093
$ws.Range("C23:I23").Value = $arr1
$arr2[0,0] = 1.930556261825757
$arr2[0,1] = 0.1406959456117818
$arr2[0,2] = 0.5378780924832256
$arr2[0,3] = 1.40964383696744
$ws.Range("K23:N23").Value = $arr2

# Row 24
$arr1[0,0] = 0.1590429706102157
$arr1[0,1] = 0.05236938097368338
$arr1[0,2] = 0.1092183336995731
$arr1[0,3] = 1.294248684692846
$arr1[0,4] = 1.167213605577899
$arr1[0,5] = 1.0995743164071
$arr1[0,6] = 1.241553048747306
$ws.Range("C24:I24").Value = $arr1
$arr2[0,0] = 1.70983044414487
$arr2[0,1] = 0.1407328895753537
$arr2[0,2] = 0.4926186894122253
$arr2[0,3] = 1.433920582487779
$ws.Range("K24:N24").Value = $arr2

# Row 25
$arr1[0,0] = 0.1555862212261303
$arr1[0,1] = 0.05187182594049489
$arr1[0,2] = 0.1084971694954042
$arr1[0,3] = 1.274655855017969
$arr1[0,4] = 1.150580113386425
$arr1[0,5] = 1.101663945160368
$arr1[0,6] = 1.231246421502377
$ws.Range("C25:I25").Value = $arr1
$arr2[0,0] = 1.472911236598691
$arr2[0,1] = 0.1410220581960431
$arr2[0,2] = 0.4444112083560015
$arr2[0,3] = 1.462274041405969
$ws.Range("K25:N25").Value = $arr2

